$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column S (year 2022) to the table, mirroring the formatting of
# column R (the previous last year column) for each row.
$ws.Range("R4:R14").Copy($ws.Range("S4:S14"))

$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 0
$ws.Range("S6").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("S8").Value = 0
$ws.Range("S9").Value = 0
$ws.Range("S10").Value = 0
$ws.Range("S11").Value = 0
$ws.Range("S12").Value = 0
$ws.Range("S13").Value = 0
$ws.Range("S14").Value = 0

# Update the active selection to match the authored state.
$ws.Range("R17").Select()
